# Add two new year columns (N = 2021, O = 2022) to the sanitation-access
# table, mirroring the formatting already used by the existing year
# columns (D:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3 — blank separator cells under the thick bottom rule, same style
# as the existing K3:M3 cells (Times New Roman 10, bottom border).
# ---------------------------------------------------------------------
foreach ($addr in @("N3", "O3")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 10
    $cell.Font.Bold = $false
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = -4138
    $cell.Borders.Item(9).ColorIndex = 1
}

# ---------------------------------------------------------------------
# Row 4 — year headers (bold Times New Roman 9, bottom border, right
# aligned) identical to the existing D4:M4 headers.
# ---------------------------------------------------------------------
$years = @{ "N4" = 2021; "O4" = 2022 }
foreach ($addr in $years.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $years[$addr]
    $cell.NumberFormat = "General"
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4152
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = -4138
    $cell.Borders.Item(9).ColorIndex = 1
}

# ---------------------------------------------------------------------
# Row 5 — bold data row (Times New Roman 10 bold, custom "0.0" format).
# ---------------------------------------------------------------------
$row5 = @{ "N5" = 40.007977647471066; "O5" = 42.620582506455563 }
foreach ($addr in $row5.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $row5[$addr]
    $cell.NumberFormat = "0.0"
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 10
    $cell.Font.Bold = $true
}

# ---------------------------------------------------------------------
# Rows 6-13 — regular data rows (Times New Roman 10, custom "0.0"
# format, no border).
# ---------------------------------------------------------------------
$rows6to13 = @{
    "N6"  = 5.7072514621689896;  "O6"  = 8.1443914479075037
    "N7"  = 8.9893229854028949;  "O7"  = 10.715961386284755
    "N8"  = 66.307512472824584;  "O8"  = 81.977461999426666
    "N9"  = 23.475213049310256;  "O9"  = 29.828871240443185
    "N10" = 9.8045372040896162;  "O10" = 9.7218425128664112
    "N11" = 9.3737779268960448;  "O11" = 8.6167819403064012
    "N12" = 70.457032471318783;  "O12" = 69.915337594090886
    "N13" = 98.411252120183207;  "O13" = 99.08571752721997
}
foreach ($addr in $rows6to13.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $rows6to13[$addr]
    $cell.NumberFormat = "0.0"
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 10
    $cell.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Row 14 — totals row (Times New Roman 10, custom "0.0" format, bottom
# border matching the rest of the row).
# ---------------------------------------------------------------------
$row14 = @{ "N14" = 63.900563564170795; "O14" = 64.805252627098838 }
foreach ($addr in $row14.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $row14[$addr]
    $cell.NumberFormat = "0.0"
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 10
    $cell.Font.Bold = $false
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = -4138
    $cell.Borders.Item(9).ColorIndex = 1
}

# ---------------------------------------------------------------------
# Move the active selection, matching the author's final cursor spot.
# ---------------------------------------------------------------------
$ws.Range("P8").Select()
